# Updated symbol list on Thu Jan 26 18:41:32 UTC 2023 with GitHub Actions
#
# All data cells on the sheet are stored as text (t="inlineStr") even though
# many of them look numeric (prices, percentages). Assigning a plain string
# like "305.46" to a Range.Value normally lets Excel's type-coercion turn it
# into a real number (or a percentage fraction for "1.55%"), which would
# change the cell's stored representation. To keep these as literal text -
# matching the source data - we briefly force the cell to Text number format
# before writing the value, then restore the "Normal" style so no stray
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $value) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - BNB
Set-TextValue "D2" "305.46"
Set-TextValue "E2" "1.55%"

# Row 3 - OKB
Set-TextValue "D3" "35.89"
Set-TextValue "E3" "1.12%"

# Row 4 - HuobiToken
Set-TextValue "E4" "-0.05%"

# Row 5 - Cronos
Set-TextValue "D5" "0.08097"
Set-TextValue "E5" "1.43%"

# Row 6 - FTXToken
Set-TextValue "D6" "1.923"
Set-TextValue "E6" "1.48%"

# Row 7 - GateToken
Set-TextValue "D7" "4.152"
Set-TextValue "E7" "2.41%"

# Row 8 - KuCoinToken
Set-TextValue "D8" "7.842"
Set-TextValue "E8" "1.12%"

# Row 9 - MXToken
Set-TextValue "D9" "0.9321"
Set-TextValue "E9" "0.41%"

# Row 10 - LiechtensteinCryptoassetsExchange
Set-TextValue "D10" "0.1264"
Set-TextValue "E10" "-13.16%"

# Row 11 - WazirX
Set-TextValue "D11" "0.1923"
Set-TextValue "E11" "0.98%"

# Row 12 - MandalaExchangeToken
Set-TextValue "D12" "0.09197"
Set-TextValue "E12" "1.88%"

# Row 13 - BitrueCoin
Set-TextValue "D13" "0.03495"
Set-TextValue "E13" "0.09%"

# Row 14 - BitMartToken
Set-TextValue "D14" "0.09928"
Set-TextValue "E14" "0.72%"

# Row 15 - BitForexToken
Set-TextValue "D15" "0.001419"
Set-TextValue "E15" "1.71%"

# Row 16 - TigerCash
Set-TextValue "D16" "0.006655"
Set-TextValue "E16" "16.33%"

# Row 17 - LEO
Set-TextValue "D17" "3.615"
Set-TextValue "E17" "2.39%"

# Row 18 - BTSEToken
Set-TextValue "D18" "3.229"
Set-TextValue "E18" "7.92%"

# Row 19 - BitpandaEcosystemToken
Set-TextValue "D19" "0.3409"
Set-TextValue "E19" "-1.11%"

# Row 20 / 21 swap places: ProBitToken <-> MCDex
Set-TextValue "B20" "MCDex"
Set-TextValue "C20" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D20" "5.175"
Set-TextValue "E20" "2.84%"

Set-TextValue "B21" "ProBitToken"
Set-TextValue "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D21" "0.1303"
Set-TextValue "E21" "-0.06%"

# Row 22 - ZBToken
Set-TextValue "D22" "0.2529"
Set-TextValue "E22" "5.43%"

# Row 23 - CoinExToken
Set-TextValue "D23" "0.04400"
Set-TextValue "E23" "-2.14%"

# Row 24 - BitKan
Set-TextValue "D24" "0.001235"
Set-TextValue "E24" "1.80%"

# Row 25 - HotbitToken
Set-TextValue "D25" "0.004734"
Set-TextValue "E25" "-0.50%"

# Row 26 - NitroEx
Set-TextValue "E26" "5.69%"

# Row 27 - UpBots
Set-TextValue "E27" "3.31%"

# Row 39 - One
Set-TextValue "D39" "0.01997"
Set-TextValue "E39" "9.59%"

# Row 40 - IDEX
Set-TextValue "D40" "0.05213"
Set-TextValue "E40" "10.03%"

# Row 41 - KickToken
Set-TextValue "D41" "0.007535"
Set-TextValue "E41" "3.21%"

# Row 42 - Dexo
Set-TextValue "D42" "0.01010"
Set-TextValue "E42" "-4.40%"

# Row 43 - BKEXToken
Set-TextValue "D43" "0.1372"
Set-TextValue "E43" "3.45%"

# Row 44 - CEJI
Set-TextValue "D44" "0.002100"
Set-TextValue "E44" "-0.47%"

# Row 45 - LocalTraders
Set-TextValue "D45" "0.01070"
Set-TextValue "E45" "-2.16%"

# Row 46 - CoinLion
Set-TextValue "D46" "0.00006326"
Set-TextValue "E46" "1.70%"

# Row 47 - Kangarootoken
Set-TextValue "E47" "-0.13%"

# Row 48 - BOLO
Set-TextValue "D48" "64.96"
Set-TextValue "E48" "0.45%"

# Row 49 - CoinbaseStockToken
Set-TextValue "D49" "0.001600"

# Row 50 - CryptobidCoin
Set-TextValue "D50" "0.00002099"
Set-TextValue "E50" "-0.13%"

# Row 51 - SpecialPowerGold
Set-TextValue "E51" "-0.13%"
